$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these price cells to remain Text (they would otherwise be auto-parsed as numbers)
$textCells = @("D5", "D6", "D7", "D8", "D19", "D20", "D21", "D22", "D23", "D25", "D27", "D29", "D33", "D34", "D36", "D38", "D39", "D40", "D42", "D43", "D44", "D46", "D47", "D49", "D50", "D51")
foreach ($c in $textCells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range("D2").Value = "64.336.37"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "3.411.65"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "570.42"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").Value = "156.54"
$ws.Range("E6").Value = "  -3.32%  "
$ws.Range("D7").Value = "0.615"
$ws.Range("E7").Value = "  +5.01%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "3.416.25"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("E10").Value = "  -2.33%  "
$ws.Range("E11").Value = "  -3.07%  "
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").Value = "3.998.73"
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("E15").Value = "  -3.87%  "
$ws.Range("E16").Value = "  -4.44%  "
$ws.Range("D17").Value = "64.384.44"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "3.448.73"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "6.34"
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("D20").Value = "13.95"
$ws.Range("E20").Value = "  -3.82%  "
$ws.Range("D21").Value = "374.67"
$ws.Range("E21").Value = "  -4.34%  "
$ws.Range("D22").Value = "7.96"
$ws.Range("E22").Value = "  -3.22%  "
$ws.Range("D23").Value = "0.550"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "72.21"
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("E26").Value = "  -4.45%  "
$ws.Range("D27").Value = "10.16"
$ws.Range("E27").Value = "  +6.85%  "
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("E30").Value = "  +2.58%  "
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").Value = "23.09"
$ws.Range("E33").Value = "  -2.53%  "
$ws.Range("D34").Value = "7.19"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("E35").Value = "  +5.40%  "
$ws.Range("D36").Value = "160.66"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").Value = "0.0758"
$ws.Range("E38").Value = "  -2.33%  "

# Row 39 and 40 swap: RenderToken <-> EnergySwap with new values
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "26.68"
$ws.Range("E39").Value = "  -3.34%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "6.73"
$ws.Range("E40").Value = "  +2.32%  "

$ws.Range("D41").Value = "2.845.44"
$ws.Range("E41").Value = "  -2.84%  "
$ws.Range("D42").Value = "4.60"
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("D43").Value = "42.73"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").Value = "26.31"
$ws.Range("E44").Value = "  +8.40%  "
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("D46").Value = "0.766"
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("D47").Value = "316.37"
$ws.Range("E47").Value = "  +6.07%  "
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("D49").Value = "0.109"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("D50").Value = "6.58"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").Value = "0.855"
$ws.Range("E51").Value = "  -2.22%  "
